# Commit: add the NA's under duplicate_image_filename
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate column E (duplicate_image_filename) with "NA" for every
# existing data row in the top stimuli table (rows 2-21).
for ($r = 2; $r -le 21; $r++) {
    $ws.Range("E$r").Value = "NA"
}
